$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark (it currently sits, empty,
#    in a blank paragraph further up the document).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldGoBack = $d.Bookmarks("_GoBack")
    $oldGoBack.Delete()
}

# ------------------------------------------------------------------
# 2. Update the certification/date line "Date : 6-Aug-2017" to
#    "Date : 7-Aug-2017".
# ------------------------------------------------------------------
$dateRng = $d.Content
$dateRng.Find.Execute("6-Aug-2017", $true, $false, $false, $false, $false, $true, 1, $false, "7-Aug-2017", 2) | Out-Null

# ------------------------------------------------------------------
# 3. Re-create "_GoBack" right after the "7" that was just typed (this
#    mirrors what Word itself does: it drops _GoBack at the point of
#    the most recent edit).
# ------------------------------------------------------------------
$foundRng = $d.Content
$foundRng.Find.Execute("7-Aug-2017", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterSeven = $d.Range($foundRng.Start + 1, $foundRng.Start + 1)
$d.Bookmarks.Add("_GoBack", $afterSeven) | Out-Null

# ------------------------------------------------------------------
# 4. The text replace above merges the previously separate "-Aug" and
#    "-2017" runs into a single run. Re-split them back apart (with a
#    throw-away bookmark that is immediately removed) so the run
#    layout matches the original document's structure.
# ------------------------------------------------------------------
$augRng = $d.Content
$augRng.Find.Execute("7-Aug", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterAug = $d.Range($augRng.End, $augRng.End)
$d.Bookmarks.Add("TempRunSplit", $afterAug) | Out-Null
$d.Bookmarks("TempRunSplit").Delete()
